$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A8").Value = "gggggg"
